$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.032.50'
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").Value = '1.924.76'
$ws.Range("E3").Value = '  +0.89%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.40'
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4587'
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3816'
$ws.Range("E8").Value = '  -0.11%  '

$ws.Range("E9").Value = '  +0.06%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9789'
$ws.Range("E10").Value = '  -0.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.63'
$ws.Range("E11").Value = '  +2.50%  '

$ws.Range("D12").Value = '1.955.14'
$ws.Range("E12").Value = '  +1.41%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.700'
$ws.Range("E13").Value = '  +0.20%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.965'
$ws.Range("E14").Value = '  -0.48%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07025'
$ws.Range("E15").Value = '  -0.37%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '84.80'
$ws.Range("E16").Value = '  +0.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009496'
$ws.Range("E18").Value = '  -0.91%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.70'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("E20").Value = '  +0.00%  '

$ws.Range("D21").Value = '29.055.51'
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.04'
$ws.Range("E23").Value = '  +0.74%  '

$ws.Range("D24").Value = '2.174.86'
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.064'
$ws.Range("E25").Value = '  -0.59%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '157.99'
$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.01'
$ws.Range("E27").Value = '  -1.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.599'
$ws.Range("E28").Value = '  -0.09%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.61'
$ws.Range("E29").Value = '  -0.25%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.832'
$ws.Range("E30").Value = '  +0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09326'
$ws.Range("E31").Value = '  +0.61%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8604'
$ws.Range("E32").Value = '  -0.13%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.094'
$ws.Range("E33").Value = '  -0.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.244'
$ws.Range("E34").Value = '  -0.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.014'
$ws.Range("E35").Value = '  -0.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05688'
$ws.Range("E36").Value = '  -0.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.149'
$ws.Range("E37").Value = '  +0.30%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.003'
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02044'
$ws.Range("E39").Value = '  +0.07%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.105'
$ws.Range("E40").Value = '  +12.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.436'
$ws.Range("E41").Value = '  -0.82%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5505'
$ws.Range("E42").Value = '  -0.64%  '

$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.345'
$ws.Range("E44").Value = '  -0.03%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002821'
$ws.Range("E45").Value = '  +8.15%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.188'
$ws.Range("E46").Value = '  +4.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5178'
$ws.Range("E47").Value = '  -0.78%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06933'
$ws.Range("E48").Value = '  +1.64%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '11.18'
$ws.Range("E49").Value = '  -1.72%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '110.62'
$ws.Range("E50").Value = '  -1.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.760'
$ws.Range("E51").Value = '  -0.94%  '
